$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "310.15"
Set-TextValue "E2" "0.56%"
Set-TextValue "D3" "39.52"
Set-TextValue "E3" "1.78%"
Set-TextValue "D4" "5.124"
Set-TextValue "E4" "0.47%"
Set-TextValue "D5" "0.08106"
Set-TextValue "E5" "-0.41%"
Set-TextValue "D6" "2.023"
Set-TextValue "E6" "3.14%"
Set-TextValue "D7" "8.165"
Set-TextValue "E7" "2.92%"
Set-TextValue "D8" "0.9294"
Set-TextValue "E8" "0.15%"
Set-TextValue "D9" "0.1421"
Set-TextValue "E9" "-0.83%"
Set-TextValue "D10" "0.1931"
Set-TextValue "E10" "-1.03%"
Set-TextValue "D11" "0.09035"
Set-TextValue "E11" "-0.81%"
Set-TextValue "D12" "0.03501"
Set-TextValue "E12" "-0.16%"
Set-TextValue "D13" "0.09822"
Set-TextValue "E13" "-0.13%"
Set-TextValue "D14" "0.001400"
Set-TextValue "E14" "-0.86%"
Set-TextValue "D15" "0.006052"
Set-TextValue "E15" "-0.24%"
Set-TextValue "D16" "3.861"
Set-TextValue "E16" "7.35%"
Set-TextValue "D17" "4.239"
Set-TextValue "E17" "1.01%"
Set-TextValue "E18" "-4.40%"
Set-TextValue "E19" "0.19%"
Set-TextValue "D20" "0.1342"
Set-TextValue "E20" "2.14%"
Set-TextValue "D21" "4.722"
Set-TextValue "E21" "-1.71%"
Set-TextValue "D22" "0.2428"
Set-TextValue "E22" "-0.24%"
Set-TextValue "D23" "0.04373"
Set-TextValue "E23" "-1.89%"
Set-TextValue "D24" "0.001229"
Set-TextValue "E24" "-0.87%"
Set-TextValue "D25" "0.004798"
Set-TextValue "E25" "-1.20%"
Set-TextValue "D26" "0.0001252"
Set-TextValue "E26" "-3.87%"
Set-TextValue "D27" "0.0004007"
Set-TextValue "E27" "-9.91%"
Set-TextValue "D39" "0.02129"
Set-TextValue "E39" "1.26%"
Set-TextValue "D40" "0.05094"
Set-TextValue "E40" "-0.02%"
Set-TextValue "D41" "0.007416"
Set-TextValue "E41" "-0.68%"
Set-TextValue "D42" "0.009771"
Set-TextValue "E42" "-3.64%"
Set-TextValue "D43" "0.1362"
Set-TextValue "E43" "-0.09%"
Set-TextValue "D44" "0.002133"
Set-TextValue "E44" "-0.49%"
Set-TextValue "D45" "0.008624"
Set-TextValue "E45" "-17.18%"
Set-TextValue "D46" "0.00006404"
Set-TextValue "E46" "3.01%"
Set-TextValue "E47" "-0.02%"
Set-TextValue "D48" "0.001002"
Set-TextValue "E48" "-37.49%"
Set-TextValue "D49" "0.002556"
Set-TextValue "E49" "-16.52%"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "-0.02%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "-0.02%"
